# Repull data, push all data, mean calculation
# Update the "dSF" (column F) values on the sheet to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -3
    5  = -3
    22 = 2
    27 = 1
    28 = 4
    29 = 3
    30 = -5
    31 = -2
    34 = -5
    36 = -2
    37 = -7
    42 = -5
    43 = 6
    44 = -8
    48 = -3
    49 = -6
    50 = -5
    55 = -4
    56 = -5
    58 = 3
    60 = -1
    61 = -2
    62 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
